$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.629.90"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.811.93"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "37.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.23%  "
$ws.Range("E9").Value = "  -3.62%  "
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0971"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "2.073.99"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "1.839.48"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "34.588.90"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "1.365.42"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -5.54%  "
$ws.Range("E40").Value = "  +6.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "81.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.44%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.940"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  -7.02%  "
